$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.012.27"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "3.483.22"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.15%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "603.45"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "143.09"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "

$ws.Range("D7").Value = "3.482.71"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.74%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "8.18"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("E11").Value = "  -4.58%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.412"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("D13").Value = "4.069.70"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("E14").Value = "  -4.25%  "

$ws.Range("E15").Value = "  -2.23%  "

$ws.Range("D16").Value = "3.482.48"
$ws.Range("E16").Value = "  +0.36%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.117"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").Value = "66.113.02"
$ws.Range("E18").Value = "  -1.12%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "10.36"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.38%  "

$ws.Range("E20").Value = "  -3.27%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.73"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.76%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "420.28"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("E23").Value = "  -2.33%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "77.44"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "3.610.83"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -4.19%  "

$ws.Range("E28").Value = "  -5.00%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.98"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.94%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("E33").Value = "  -6.95%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "25.16"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "3.477.88"
$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("E37").Value = "  -4.17%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "5.57"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -5.52%  "

$ws.Range("E39").Value = "  -2.30%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "170.07"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0863"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.890"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("E44").Value = "  -5.43%  "

$ws.Range("E45").Value = "  -7.43%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "45.03"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "25.97"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -9.83%  "

$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("E50").Value = "  -4.23%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.930"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -3.78%  "
